# Generate Report for Handoff
# Adds a new row (row 3) to the "Overview", "zh-cn" and "de-de" sheets for the
# newly processed file "c283cab1-041b-4ad9-9040-d5d057c53b21...md", mirroring
# the existing row 2 that was created for "67d6b785-2334-4e29-a22e-57cb38d32e5d...md".

$wb = $excel.ActiveWorkbook

$repoBlobBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/833fdf6634fe671f37351d91f6c2252832707140/e2e/"

$newMdName      = "c283cab1-041b-4ad9-9040-d5d057c53b21ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdDisplay   = "e2e\c283cab1-041b-4ad9-9040-d5d057c53b21ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newStatus      = "Ready for handoff"
$overviewDate   = "2016-08-24 14:30:29"

$newZhXlf       = "c283cab1-041b-4ad9-9040-d5d057c53b21oooooooooooooooooooooooooooooooooooooooo.1285a781b62e8ffdc3dcec12c24e92f7ea5d90b4.zh-cn.xlf"
$newZhDate      = "2016-08-24 14:30:03"

$newDeXlf       = "c283cab1-041b-4ad9-9040-d5d057c53b21oooooooooooooooooooooooooooooooooooooooo.1285a781b62e8ffdc3dcec12c24e92f7ea5d90b4.de-de.xlf"
$newDeDate      = "2016-08-24 14:30:29"

# ---------------------------------------------------------------------------
# Sheet "Overview" (table3 / Overview): columns A..G
# A=File Name, B=Path And Name, C=Extension, D=Publish URL, E=zh-cn, F=de-de,
# G=Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newMdName
$wsOverview.Range("B3").Value = $newMdDisplay
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $overviewDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), ($repoBlobBase + $newMdName), $null, $null, $newMdDisplay)

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (table1): columns A..P
# A=Source File Name, B=File Extension, C=Status, D=Source Path, E=Priority,
# F=Content Duplicate, G=Latest Handoff File, H=Latest Handoff Datetime,
# I=Latest Target File, J=Latest Handback File, K=Latest Handback DateTime,
# L=Reference Tokens, M=To be localized, N=Dependency From, O=Has metadata,
# P=Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = $newMdName
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $newZhXlf
$wsZhCn.Range("H3").Value = $newZhDate
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), ($repoBlobBase + $newMdName), $null, $null, $newMdDisplay)

# ---------------------------------------------------------------------------
# Sheet "de-de" (table2): same columns as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = $newMdName
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $newDeXlf
$wsDeDe.Range("H3").Value = $newDeDate
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), ($repoBlobBase + $newMdName), $null, $null, $newMdDisplay)

# ---------------------------------------------------------------------------
# Styling to mirror row 2: hyperlink-styled filename cell + date-formatted cells
# ---------------------------------------------------------------------------
$wsOverview.Range("B3").Style = "Hyperlink"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("A3").Style = "Hyperlink"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("A3").Style = "Hyperlink"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Extend the tables so the new row is included, and autofit the Status
# columns which grew wider once "Ready for handoff" was introduced.
# ---------------------------------------------------------------------------
$wb.ActiveSheet.ListObjects | Out-Null

$tblOverview = $wsOverview.ListObjects.Item("Overview")
$tblOverview.Resize($wsOverview.Range("A1:G3"))

$tblZhCn = $wsZhCn.ListObjects.Item("zh_cn")
$tblZhCn.Resize($wsZhCn.Range("A1:P3"))

$tblDeDe = $wsDeDe.ListObjects.Item("de_de")
$tblDeDe.Resize($wsDeDe.Range("A1:P3"))

$wsOverview.Range("E1").EntireColumn.AutoFit()
$wsOverview.Range("F1").EntireColumn.AutoFit()
$wsZhCn.Range("C1").EntireColumn.AutoFit()
$wsDeDe.Range("C1").EntireColumn.AutoFit()
